# Auto-generated edit script: update numeric values in result_data_RandomForest sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = -12.6864
$ws.Cells.Item(2, 5).Value = 16.28760000000001
$ws.Cells.Item(3, 1).Value = -21.7341
$ws.Cells.Item(3, 4).Value = -7.468699999999995
$ws.Cells.Item(6, 5).Value = 16.5319
$ws.Cells.Item(12, 4).Value = -7.455299999999998
$ws.Cells.Item(14, 1).Value = -21.64980000000001
$ws.Cells.Item(16, 1).Value = -21.96399999999998
$ws.Cells.Item(18, 3).Value = -12.73570000000001
$ws.Cells.Item(19, 5).Value = 16.2737
$ws.Cells.Item(21, 1).Value = -20.17909999999997
$ws.Cells.Item(23, 1).Value = -20.17709999999998
$ws.Cells.Item(24, 3).Value = -12.2734
$ws.Cells.Item(24, 4).Value = -8.028699999999997
$ws.Cells.Item(24, 5).Value = 16.74100000000001
$ws.Cells.Item(25, 1).Value = -22.48499999999999
$ws.Cells.Item(25, 3).Value = -13.723
$ws.Cells.Item(25, 4).Value = -9.210799999999997
$ws.Cells.Item(26, 1).Value = -21.21489999999997
$ws.Cells.Item(27, 3).Value = -12.70549999999999
$ws.Cells.Item(27, 5).Value = 16.60519999999999
$ws.Cells.Item(29, 1).Value = -20.93989999999997
$ws.Cells.Item(30, 3).Value = -12.9696
$ws.Cells.Item(30, 5).Value = 15.94270000000001
$ws.Cells.Item(31, 3).Value = -12.58
$ws.Cells.Item(31, 5).Value = 16.22200000000002
$ws.Cells.Item(33, 5).Value = 16.91770000000001
$ws.Cells.Item(39, 3).Value = -12.76750000000001
$ws.Cells.Item(40, 1).Value = -19.99539999999999
$ws.Cells.Item(41, 4).Value = -8.260299999999996
$ws.Cells.Item(42, 3).Value = -12.65459999999999
$ws.Cells.Item(42, 5).Value = 16.3727
$ws.Cells.Item(48, 3).Value = -11.86059999999999
$ws.Cells.Item(50, 4).Value = -8.008500000000005
$ws.Cells.Item(51, 3).Value = -11.6871
$ws.Cells.Item(52, 3).Value = -11.3526
$ws.Cells.Item(53, 1).Value = -22.56680000000001
$ws.Cells.Item(53, 4).Value = -6.304300000000001
$ws.Cells.Item(55, 3).Value = -13.65869999999999
$ws.Cells.Item(55, 5).Value = 16.41030000000001
$ws.Cells.Item(56, 3).Value = -11.46660000000001
$ws.Cells.Item(56, 4).Value = -8.097700000000005
$ws.Cells.Item(57, 1).Value = -22.14419999999999
$ws.Cells.Item(57, 3).Value = -13.30049999999999
$ws.Cells.Item(57, 4).Value = -8.864000000000001
$ws.Cells.Item(58, 4).Value = -8.271800000000008
$ws.Cells.Item(58, 5).Value = 16.13760000000002
$ws.Cells.Item(59, 1).Value = -22.48669999999999
$ws.Cells.Item(60, 3).Value = -12.98019999999999
$ws.Cells.Item(61, 4).Value = -7.9884
$ws.Cells.Item(63, 4).Value = -7.898700000000002
$ws.Cells.Item(64, 4).Value = -7.7751
$ws.Cells.Item(65, 1).Value = -21.84349999999998
$ws.Cells.Item(65, 5).Value = 17.06210000000002
$ws.Cells.Item(69, 1).Value = -21.58269999999998
$ws.Cells.Item(70, 4).Value = -8.062800000000005
$ws.Cells.Item(70, 5).Value = 16.82869999999999
$ws.Cells.Item(72, 4).Value = -7.358599999999995
$ws.Cells.Item(73, 3).Value = -12.5151
$ws.Cells.Item(74, 3).Value = -12.66670000000001
$ws.Cells.Item(74, 5).Value = 16.76329999999999
$ws.Cells.Item(75, 5).Value = 16.39120000000001
$ws.Cells.Item(79, 1).Value = -20.68010000000001
$ws.Cells.Item(83, 1).Value = -21.77319999999999
$ws.Cells.Item(83, 5).Value = 16.7065
$ws.Cells.Item(84, 5).Value = 16.8437
$ws.Cells.Item(86, 4).Value = -8.528900000000007
$ws.Cells.Item(86, 5).Value = 16.0086
$ws.Cells.Item(89, 3).Value = -10.3358
$ws.Cells.Item(89, 4).Value = -5.612300000000002
$ws.Cells.Item(90, 3).Value = -12.2266
$ws.Cells.Item(91, 1).Value = -21.36350000000003
$ws.Cells.Item(92, 3).Value = -10.6232
$ws.Cells.Item(93, 1).Value = -20.93869999999998
$ws.Cells.Item(96, 5).Value = 16.0586
$ws.Cells.Item(97, 5).Value = 16.90510000000002
$ws.Cells.Item(98, 4).Value = -8.772499999999994
$ws.Cells.Item(100, 1).Value = -22.3749
$ws.Cells.Item(100, 4).Value = -8.537600000000003
$ws.Cells.Item(102, 4).Value = -7.350499999999994

$wb.Save()
